$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# Grading checklist: mark the "enemies" and "power-ups" requirements (rows 29,
# 31, 40, 42) as fully implemented now that they're in the build. Each of
# these cells is a data-validation dropdown of "No Credit" / "Partial
# Credit" / "Full Credit"; the A-column score formulas recalc automatically.
$ws.Range("B29").Value = "Full Credit"
$ws.Range("B31").Value = "Full Credit"
$ws.Range("B40").Value = "Full Credit"
$ws.Range("B42").Value = "Full Credit"

# Move the viewport/selection to where the author was working.
$ws.Activate()
$ws.Range("C44").Select()
